# Weekly update: insert a new data row for "Perejil" (Terminal La Palmera de
# La Serena) ahead of the existing row 213, pushing the rest of the table
# down by one row (old row 213 becomes 214, ..., old row 256 becomes 257).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 213; Excel shifts rows 213:256 down to 214:257
# and copies formatting from the row above, matching the row style pattern
# used throughout this table.
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with the new weekly data point.
$ws.Cells.Item(213, 1).Value = 8
$ws.Cells.Item(213, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = 45211
$ws.Cells.Item(213, 4).NumberFormat = $ws.Cells.Item(214, 4).NumberFormat
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 6).Value = 100112044
$ws.Cells.Item(213, 7).Value = "Perejil"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Primera"
$ws.Cells.Item(213, 10).Value = 2400
$ws.Cells.Item(213, 11).Value = 1500
$ws.Cells.Item(213, 12).Value = 2000
$ws.Cells.Item(213, 13).Value = 1750
$ws.Cells.Item(213, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(213, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(213, 16).Value = 1167
$ws.Cells.Item(213, 17).Value = 1.5
$ws.Cells.Item(213, 18).Value = "Hortaliza"
